# The author selected the whole of column B (header + data) and ran
# Find & Replace, swapping "." for "," (a locale-style decimal-separator
# tweak). Excel's Replace rewrites any cell whose displayed text contains
# the search string into a literal string cell - so:
#   - the rich-text header in B1 ("Durchschnittl. Mietpreis pro m<sup>2</sup> 3-Zimmer")
#     loses its run formatting and becomes the plain string
#     "Durchschnittl, Mietpreis pro m2 3-Zimmer"
#   - every data cell in column B whose number displays with a decimal
#     point (e.g. 15.4) turns into the text "15,4"
#   - whole numbers (e.g. 17, 20, 14) have no "." in their displayed text,
#     so Replace leaves them as untouched numeric cells
# This matches the resulting workbook exactly, including which new shared
# strings get created and in what order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Range("B1:B1048576")
[void]$col.Select()
[void]$col.Replace(".", ",")
